# Apply the diff described in the task:
#  - G2: 2 -> 1
#  - H2: "无法获取有效的视频下载地址 (可能是Blob加密流)" -> "" (empty string)
#  - Append six new data rows (17-22)
#  - Dimension / ignoredErrors range grows to A1:I22 (handled automatically by the engine)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- existing-row edits -----------------------------------------------
$ws.Range("G2").Value = 1

# Writing Value = "" clears the cell entirely (becomes a true blank cell,
# not an empty-string text cell), which doesn't match the source (a
# literal empty-string text cell). The "'" (lone quote-prefix) trick
# forces Excel to store an empty *text* value; resetting the style back
# to Normal afterwards drops the quote-prefix formatting flag it leaves
# behind so the cell is indistinguishable from any other empty-string
# text cell in the sheet.
$ws.Range("H2").Value = "'"
$ws.Range("H2").Style = "Normal"

# --- new rows -----------------------------------------------------------
$newRows = @(
    @(16, "扣子2.0全新升级！快来看看我开发的网页应用吧", "视频", "#扣子编程", "https://www.toutiao.com/video/7597347992571478528/", "2026-01-20", 1, "videos\2026-01-20\16.mp4", "已提取"),
    @(17, "Claude Cowork（实测）：这新功能太上头了，我可能要玩上瘾！", "视频", "#ClaudeCowork", "https://www.toutiao.com/video/7595140659170771491/", "2026-01-20", 1, "videos\2026-01-20\17.mp4", "已提取"),
    @(18, "用上这个Skill，你的Claude Code/Codex 将会比别人快5倍 -- 用分布式思维驯服AI任务编排", "文章", "无分类", "https://www.toutiao.com/article/7594654759382991398/", "2026-01-20", 1, "articles\2026-01-20\18.txt", ""),
    @(19, "为什么Agent总是Demo猛如龙实战一条虫？", "文章", "无分类", "https://www.toutiao.com/article/7586626602591732262/", "2026-01-20", 1, "articles\2026-01-20\19.txt", ""),
    @(20, "手把手教你用上开源版Claude Code，人人都可以体验编程Agent的魅力了。", "文章", "无分类", "https://www.toutiao.com/article/7594454191276294719/", "2026-01-20", 1, "articles\2026-01-20\20.txt", ""),
    @(21, "结合 Cluade Code，我为项目设计了一个 SKILL，AI 输出代码可用率达到了 99%，谈谈真实的使用感受", "文章", "无分类", "https://www.toutiao.com/article/7596202303636439552/", "2026-01-20", 1, "articles\2026-01-20\21.txt", "")
)

$startRow = 17
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $data = $newRows[$i]

    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
    $ws.Cells.Item($r, 5).Value = $data[4]
    # column F holds the date as literal text ("2026-01-20"), not a real
    # date value. A bare Value assignment lets COM auto-convert the
    # recognisable date string to a date serial, so instead prefix it
    # with a quote (forces text entry) and then drop the resulting
    # quote-prefix style flag, leaving a plain unstyled text cell - same
    # trick used for the emptied-out H2 cell above.
    $ws.Cells.Item($r, 6).Value = "'" + $data[5]
    $ws.Cells.Item($r, 6).Style = "Normal"
    $ws.Cells.Item($r, 7).Value = $data[6]
    $ws.Cells.Item($r, 8).Value = $data[7]

    # column I ("audio status") is blank text for article rows in the
    # source. Same empty-text trick as H2 above.
    if ($data[8] -eq "") {
        $ws.Cells.Item($r, 9).Value = "'"
        $ws.Cells.Item($r, 9).Style = "Normal"
    } else {
        $ws.Cells.Item($r, 9).Value = $data[8]
    }
}
